$d = $word.ActiveDocument

# Remove the first 9 paragraphs (title/abstract block + trailing page-break
# paragraph) that precede the main "Title: Predicting NBA Success with
# Combined Psychological and Physical Statistics" paragraph.
$start = $d.Paragraphs.Item(1).Range.Start
$end = $d.Paragraphs.Item(9).Range.End
$r = $d.Range($start, $end)
$r.Delete()
